# Update "countries & provincias Spain" data in the Pais sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 22:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 333017
$ws.Range("C4").Value = 21660
$ws.Range("E4").Value = 306471

# Row 5 - Espana
$ws.Range("B5").Value = 130854
$ws.Range("C5").Value = 4686
$ws.Range("E5").Value = 80256
$ws.Range("G5").Value = 571
$ws.Range("H5").Value = 12518

# Row 7 - Alemania
$ws.Range("B7").Value = 100024
$ws.Range("C7").Value = 3932
$ws.Range("E7").Value = 69748
$ws.Range("G7").Value = 132
$ws.Range("H7").Value = 1576

# Row 8 - Francia
$ws.Range("B8").Value = 92839
$ws.Range("C8").Value = 2886
$ws.Range("D8").Value = 16183
$ws.Range("E8").Value = 68578
$ws.Range("G8").Value = 518
$ws.Range("H8").Value = 8078

# Row 19 - Brasil
$ws.Range("B19").Value = 11130
$ws.Range("C19").Value = 770
$ws.Range("E19").Value = 10517
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = 486

# Row 47 - Emiratos Arabes Unidos
$ws.Range("B47").Value = 1799
$ws.Range("C47").Value = 294
$ws.Range("E47").Value = 1645
